# wrapping up test file audit
#
# The "optimization_parameters" sheet had a stray leftover row (a label
# "Sheet" with the values 3 / 4 next to it) that doesn't belong with the
# rest of the optimization parameters. Remove it; everything below moves
# up one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

$ws.Range("A16:C16").EntireRow.Delete() | Out-Null

# Leave the selection where the removed row used to be, same as Excel
# does after an in-sheet row delete.
$ws.Rows.Item(16).Select() | Out-Null
